$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.660.02"
$ws.Range("E2").Value = "  -7.00%  "

$ws.Range("D3").Value = "'1.695.79"
$ws.Range("E3").Value = "  -5.69%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "'219.89"
$ws.Range("E5").Value = "  -5.03%  "

$ws.Range("D6").Value = "'0.5121"
$ws.Range("E6").Value = "  -13.06%  "

$ws.Range("D8").Value = "'0.2641"
$ws.Range("E8").Value = "  -4.62%  "

$ws.Range("D9").Value = "'22.16"
$ws.Range("E9").Value = "  -4.67%  "

$ws.Range("D10").Value = "'0.06295"
$ws.Range("E10").Value = "  -7.48%  "

$ws.Range("D11").Value = "'0.07343"
$ws.Range("E11").Value = "  -2.36%  "

$ws.Range("D12").Value = "'1.701.19"
$ws.Range("E12").Value = "  -5.32%  "

$ws.Range("D13").Value = "'4.523"
$ws.Range("E13").Value = "  -5.47%  "

$ws.Range("D14").Value = "'0.5787"
$ws.Range("E14").Value = "  -6.64%  "

$ws.Range("D15").Value = "'1.927.42"
$ws.Range("E15").Value = "  -5.63%  "

$ws.Range("D16").Value = "'0.000008447"
$ws.Range("E16").Value = "  -7.28%  "

$ws.Range("D17").Value = "'65.56"
$ws.Range("E17").Value = "  -13.20%  "

$ws.Range("D18").Value = "'26.680.88"
$ws.Range("E18").Value = "  -6.86%  "

$ws.Range("D19").Value = "'5.000"
$ws.Range("E19").Value = "  -8.82%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  -4.50%  "

$ws.Range("D22").Value = "'186.74"
$ws.Range("E22").Value = "  -11.58%  "

$ws.Range("D23").Value = "'6.248"
$ws.Range("E23").Value = "  -8.48%  "

$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").Value = "'144.81"
$ws.Range("E25").Value = "  -5.74%  "

$ws.Range("D26").Value = "'7.516"
$ws.Range("E26").Value = "  -5.77%  "

$ws.Range("D27").Value = "'0.1157"
$ws.Range("E27").Value = "  -8.76%  "

$ws.Range("D28").Value = "'15.80"
$ws.Range("E28").Value = "  -3.86%  "

$ws.Range("D29").Value = "'1.350"
$ws.Range("E29").Value = "  -5.16%  "

$ws.Range("D30").Value = "'0.05657"
$ws.Range("E30").Value = "  -7.53%  "

$ws.Range("D31").Value = "'1.341"
$ws.Range("E31").Value = "  -5.92%  "

$ws.Range("E32").Value = "  -7.28%  "

$ws.Range("D33").Value = "'3.495"
$ws.Range("E33").Value = "  -8.80%  "

$ws.Range("D34").Value = "'1.650"
$ws.Range("E34").Value = "  -4.91%  "

$ws.Range("D35").Value = "'1.021"
$ws.Range("E35").Value = "  -3.12%  "

$ws.Range("D36").Value = "'0.6002"
$ws.Range("E36").Value = "  -6.65%  "

$ws.Range("D37").Value = "'2.359"
$ws.Range("E37").Value = "  -5.67%  "

$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("D39").Value = "'1.103.60"
$ws.Range("E39").Value = "  -3.62%  "

$ws.Range("D40").Value = "'0.01612"
$ws.Range("E40").Value = "  -5.04%  "

$ws.Range("D41").Value = "'0.8602"
$ws.Range("E41").Value = "  -3.03%  "

$ws.Range("D42").Value = "'5.843"
$ws.Range("E42").Value = "  -10.34%  "

$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Value = "'99.74"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("E45").Value = "  -5.03%  "

$ws.Range("D46").Value = "'0.00000000114"
$ws.Range("E46").Value = "  +2.59%  "

$ws.Range("D47").Value = "'56.62"
$ws.Range("E47").Value = "  -5.99%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").Value = "'8.111"
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").Value = "'0.05241"
$ws.Range("E50").Value = "  -4.01%  "

$ws.Range("D51").Value = "'0.4323"
$ws.Range("E51").Value = "  -3.51%  "
